$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 (the last existing data row) is copied down to the new row 20,
# carrying over all the static columns (A,B,C,E-P,R) and formatting.
$ws.Range("A19:T19").Copy($ws.Range("A20:T20"))

# Now update the per-row varying fields (D = Fecha, Q = Unidad de
# comercialización, S = Precio $/Kg, T = Kg / unidad) for rows 17-20 to
# their new target values.

# Row 17: date + "bandeja" unit
$ws.Range("D17").Value = 45029
$ws.Range("Q17").Value = "$/bandeja 18 kilos granel"
$ws.Range("S17").Value = 528
$ws.Range("T17").Value = 18

# Row 18: date + "caja 15 kilos empedrada" unit
$ws.Range("D18").Value = 44363
$ws.Range("Q18").Value = "$/caja 15 kilos empedrada"
$ws.Range("S18").Value = 633
$ws.Range("T18").Value = 15

# Row 19: date + "bandeja" unit
$ws.Range("D19").Value = 45027
$ws.Range("Q19").Value = "$/bandeja 18 kilos granel"
$ws.Range("S19").Value = 528
$ws.Range("T19").Value = 18

# Row 20 (new row): date + "caja 18 kilos granel" unit
$ws.Range("D20").Value = 44316
$ws.Range("Q20").Value = "$/caja 18 kilos granel"
$ws.Range("S20").Value = 528
$ws.Range("T20").Value = 18
